$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for rows 2..22 (columns A..F), derived from the
# "Redo report after Jan" report rebuild. Row index in $data corresponds
# to worksheet row (row r = index + 2).
$data = @(
    @(0, 4, "e5de31f2-fc40-4062-9135-3aba73149106", "4DabfkFLBojKxKol8UQf", 4, "3K5o2d8SL2y20FQOPgFK"),
    @(1, 3, "e5de31f2-fc40-4062-9135-3aba73149106", "1njZxACnAGt6T5qLMRlL", 2, "5dXobnKMzfF8oKkUeWCc"),
    @(2, 4, "e5de31f2-fc40-4062-9135-3aba73149106", "3BPbSeLlA1pP4U3oOUfy", 4, "6GXbN8XhMCIk1jY6jIYF"),
    @(3, 3, "e5de31f2-fc40-4062-9135-3aba73149106", "0YwyBZxWsQpHsGv8LS61", 2, "7KGaWT2rYsGwfUQ9Q5ET"),
    @(4, 2, "a7037b14-2e6d-4f7d-bac8-ed02d7df5329", "6VfrcjCKMIYNp2ZeSZNs", 1, "7VIAVgkaOiAjQD5UsCNe"),
    @(5, 4, "19b57dbb-8b65-48f2-a217-3287973988cd", "4F14IrWf2S6Uwzm7ajl4", 5, "C49JBc2Um2B1khgodGze"),
    @(6, 3, "e2e13763-bad1-40d3-9c2c-52c7e9d12d9f", "5C5kkshRarmXiRx92JRx", 5, "GB6L3zRFWhcORe0n1R9D"),
    @(7, 4, "19b57dbb-8b65-48f2-a217-3287973988cd", "5QLfEePCnZ99zOevCDGL", 3, "L9p36CvLC15ACnwDUgx4"),
    @(8, 2, "e2e13763-bad1-40d3-9c2c-52c7e9d12d9f", "1SUpmaSAQQNMdD8FibPB", 4, "MDGLCGGLsJZCqma8nnGv"),
    @(9, 3, "19b57dbb-8b65-48f2-a217-3287973988cd", "29sPSf4gbZUmlbkLou5s", 4, "MIa0jFtbSCU2El3EFDZu"),
    @(10, 4, "a7037b14-2e6d-4f7d-bac8-ed02d7df5329", "6gwgrl1Y6W32aOHhxbVS", 2, "N65inbIwJ0l9IefyIyPz"),
    @(11, 4, "a7037b14-2e6d-4f7d-bac8-ed02d7df5329", "0uSshiYpzXHdIn1okZsk", 4, "QIHJ1odwoylsYRzFbHo7"),
    @(12, 5, "e2e13763-bad1-40d3-9c2c-52c7e9d12d9f", "2kRptcQx5PSo54pmymCi", 3, "VMsKz72kCWNGJP6ZGlZ5"),
    @(13, 4, "e2e13763-bad1-40d3-9c2c-52c7e9d12d9f", "0Qbwb7Lis9pOLUVOSIrO", 3, "VaTsU5q2TjosMgFloSxj"),
    @(14, 4, "19b57dbb-8b65-48f2-a217-3287973988cd", "0vY7sjGJiyicoohSiEQM", 3, "VxZKHXdW6dwCXEcHiAod"),
    @(15, 3, "19b57dbb-8b65-48f2-a217-3287973988cd", "4o2MS9PcRonKh1JZMUXr", 2, "Z5fokkWGJnEX02EjVlLC"),
    @(16, 3, "a7037b14-2e6d-4f7d-bac8-ed02d7df5329", "3oCaNI7N11UjdwsgwiY6", 3, "oaxzMeJDa04t3NlPnZzJ"),
    @(17, 2, "a7037b14-2e6d-4f7d-bac8-ed02d7df5329", "9cfqp5ZF6aZKsFMPR82A", 3, "phxct5YusNvHcQksdyOb"),
    @(18, 2, "a7037b14-2e6d-4f7d-bac8-ed02d7df5329", "1DMpSYJYgWlOR38X57hm", 3, "q2CAvgjDljuJAVBSUdMW"),
    @(19, 5, "e5de31f2-fc40-4062-9135-3aba73149106", "0bIsWmiOeRXEyKdMrHB1", 5, "qM69gpRHSppZgpOisLcr"),
    @(20, 4, "e2e13763-bad1-40d3-9c2c-52c7e9d12d9f", "2yLwlpCr3EttnrWODivG", 3, "zduB5CFox0wV2nvJUVeq")
)

$firstDataRow = 2
$lastExistingRow = 11
$lastNewRow = $firstDataRow + $data.Count - 1

# Extend formatting (border/font/alignment) of column A down to the new
# rows by copying the format from A11 (which already carries the bordered
# "index" style) before writing values into the newly added rows.
if ($lastNewRow -gt $lastExistingRow) {
    $ws.Range("A$lastExistingRow").Copy() | Out-Null
    $destA = $ws.Range("A" + ($lastExistingRow + 1) + ":A" + $lastNewRow)
    $destA.PasteSpecial(-4122) | Out-Null
}

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $firstDataRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
}
